$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.293.31"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.865.80"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'0.7091"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").Value = "'237.81"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.07856"
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("D9").Value = "'0.3060"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "'25.05"
$ws.Range("E10").Value = "  +7.20%  "
$ws.Range("D11").Value = "'0.08164"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.883.22"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'5.221"
$ws.Range("D14").Value = "'0.7175"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'89.16"
$ws.Range("D16").Value = "29.328.37"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'5.805"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'241.18"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'0.000007791"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.109.98"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'7.553"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "'162.15"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'8.919"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").Value = "'0.1443"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "'1.911"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").Value = "'1.369"
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("D31").Value = "'1.476"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'4.304"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("D33").Value = "'4.041"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "'0.05191"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "'0.7157"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "'1.005"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'0.01846"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").Value = "'2.694"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "1.168.85"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "'0.9142"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "'5.996"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "'71.24"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'0.4265"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "'0.9995"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'101.99"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "'0.5353"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").Value = "'1.744"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "'9.184"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "'6.977"
$ws.Range("E51").Value = "  -0.21%  "

Write-Output "Applied 99 cell updates"
